$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("add_cart")
$ws.Range("O2").Value = "311,60 €"
$ws.Range("Q2").Value = "311,60 €"
$ws.Range("N6").Select()
